# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
#
# Price (column D) cells in this sheet are always stored as text (e.g. some
# prices use "."-as-thousands-separator like "69.457.36"), so before writing
# any numeric-looking price we force the cell to Text format ("@") first -
# otherwise Excel would silently reinterpret a value like "191.58" as a
# floating point number. Percent cells (column E) already contain a "%" sign
# plus padding spaces, so they are never misread as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.457.36'
$ws.Range('E2').Value = '  +2.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.384.63'
$ws.Range('E3').Value = '  +4.63%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '191.58'
$ws.Range('E5').Value = '  +4.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '593.13'
$ws.Range('E6').Value = '  +2.45%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.606'
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.135'
$ws.Range('E9').Value = '  +2.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.79'
$ws.Range('E10').Value = '  +3.09%  '
$ws.Range('E11').Value = '  +2.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.974.59'
$ws.Range('E12').Value = '  +4.76%  '
$ws.Range('E13').Value = '  +1.18%  '
$ws.Range('E14').Value = '  +3.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '69.499.49'
$ws.Range('E16').Value = '  +1.98%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.388.14'
$ws.Range('E17').Value = '  +5.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '449.87'
$ws.Range('E18').Value = '  +13.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.84'
$ws.Range('E19').Value = '  +1.60%  '
$ws.Range('E20').Value = '  +2.65%  '
$ws.Range('E21').Value = '  +3.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.40'
$ws.Range('E22').Value = '  +5.94%  '
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.525.20'
$ws.Range('E24').Value = '  +4.49%  '

# Rows 25 and 26 swap contents: PEPE (previously row 26) moves to row 25 and
# Polygon (previously row 25) moves to row 26; the row index in column A
# stays put for both rows.
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000123'
$ws.Range('E25').Value = '  +4.11%  '
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.522'
$ws.Range('E26').Value = '  +1.31%  '
$ws.Range('E27').Value = '  +2.06%  '
$ws.Range('E28').Value = '  -1.30%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  +1.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '23.42'
$ws.Range('E31').Value = '  +3.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.65'
$ws.Range('E32').Value = '  +1.67%  '
$ws.Range('E33').Value = '  +3.10%  '
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  +5.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '165.29'
$ws.Range('E37').Value = '  +2.80%  '
$ws.Range('E38').Value = '  +3.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.54'
$ws.Range('E39').Value = '  +4.05%  '
$ws.Range('E40').Value = '  +2.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.60'
$ws.Range('E41').Value = '  +1.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.58'
$ws.Range('E42').Value = '  +1.66%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.754.10'
$ws.Range('E43').Value = '  +5.38%  '
$ws.Range('E44').Value = '  +2.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '25.56'
$ws.Range('E45').Value = '  +3.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0691'
$ws.Range('E46').Value = '  +0.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.84'
$ws.Range('E47').Value = '  +0.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '341.30'
$ws.Range('E48').Value = '  +2.23%  '
$ws.Range('E49').Value = '  +2.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.02'
$ws.Range('E50').Value = '  +7.84%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.02'
$ws.Range('E51').Value = '  +5.96%  '
